$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 15 (shifts existing row 15 and below down by one),
# representing the new daily sales record for day 14 of May/2025.
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 36409.18
$ws.Cells.Item(15, 3).Value = 5
$ws.Cells.Item(15, 4).Value = 2025
$ws.Cells.Item(15, 5).Value = "05/2025"
